$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(110)
$rng = $p.Range
$found = $rng.Find.Execute("1hr", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$one = $d.Range($rng.Start, $rng.Start + 1)
$one.Text = "2"
